# Apply the price/volume refresh for cryptos.xlsx (GitHub Actions update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.329.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "'1.856.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.57%  "
$ws.Range("D5").Value = "'314.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").Value = "'0.4621"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").Value = "'0.3699"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("D9").Value = "'0.07323"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("D10").Value = "'0.8825"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.22%  "
$ws.Range("D11").Value = "'0.07824"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "'19.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "'1.873.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "'5.383"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").Value = "'6.532"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "'91.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D18").Value = "'0.000008869"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.13%  "
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("E20").Value = "  +1.81%  "
$ws.Range("D21").Value = "'27.355.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").Value = "'5.113"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("D24").Value = "'2.082.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "'1.888"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.70%  "
$ws.Range("D26").Value = "'152.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").Value = "'2.070"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("D29").Value = "'5.118"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "'116.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("D31").Value = "'0.08849"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "'0.7663"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.27%  "
$ws.Range("D33").Value = "'2.997"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("D34").Value = "'1.169"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.07%  "
$ws.Range("D35").Value = "'4.486"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("D36").Value = "'2.620"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.97%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.080"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01963"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("E39").Value = "  +1.88%  "
$ws.Range("D40").Value = "'0.05216"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("D41").Value = "'7.028"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.84%  "
$ws.Range("D42").Value = "'0.5148"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").Value = "'0.1638"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("D45").Value = "'0.4828"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").Value = "'10.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("D48").Value = "'103.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("D51").Value = "'65.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.01%  "
